$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill column C (Checks) with "ok" for every criteria row except the header
# and the row that already has it (row 8), matching the final reviewed state.
foreach ($r in 2..23) {
    if ($r -ne 8) {
        $ws.Cells.Item($r, 3).Value = "ok"
    }
}

# Leave the final selection on C4, matching the saved state of the workbook.
$ws.Range("C4").Select()
